$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.397.51"
$ws.Range("E2").Value = "  +0.26%  "
$ws.Range("D3").Value = "1.843.29"
$ws.Range("E3").Value = "  +0.11%  "
$ws.Range("D4").Value = "0.9988"
$ws.Range("E4").Value = "  -0.18%  "
$ws.Range("D5").Value = "240.13"
$ws.Range("E5").Value = "  -0.08%  "
$ws.Range("D6").Value = "0.6331"
$ws.Range("E6").Value = "  +0.93%  "
$ws.Range("D7").Value = "1.000"
$ws.Range("E7").Value = "  -0.09%  "
$ws.Range("D8").Value = "0.07475"
$ws.Range("E8").Value = "  +0.05%  "
$ws.Range("D9").Value = "0.2906"
$ws.Range("E9").Value = "  +0.46%  "
$ws.Range("E10").Value = "  +3.00%  "
$ws.Range("D11").Value = "0.07745"
$ws.Range("E11").Value = "  +0.40%  "
$ws.Range("D12").Value = "1.847.97"
$ws.Range("E12").Value = "  +0.40%  "
$ws.Range("D13").Value = "4.987"
$ws.Range("E13").Value = "  +0.06%  "
$ws.Range("D14").Value = "0.6790"
$ws.Range("E14").Value = "  +0.29%  "
$ws.Range("D15").Value = "0.00001023"
$ws.Range("E15").Value = "  -0.34%  "
$ws.Range("D16").Value = "82.06"
$ws.Range("E16").Value = "  +0.04%  "
$ws.Range("D17").Value = "6.267"
$ws.Range("E17").Value = "  +2.54%  "
$ws.Range("D18").Value = "29.410.03"
$ws.Range("E18").Value = "  +0.15%  "
$ws.Range("D19").Value = "230.00"
$ws.Range("E19").Value = "  +0.72%  "
$ws.Range("E20").Value = "  +0.72%  "
$ws.Range("D21").Value = "0.9997"
$ws.Range("E21").Value = "  -0.15%  "
$ws.Range("D22").Value = "7.425"
$ws.Range("E22").Value = "  +0.73%  "
$ws.Range("E23").Value = "  -0.20%  "
$ws.Range("D24").Value = "158.18"
$ws.Range("E24").Value = "  -0.45%  "
$ws.Range("D25").Value = "8.500"
$ws.Range("E25").Value = "  +1.49%  "
$ws.Range("E26").Value = "  -1.73%  "
$ws.Range("D27").Value = "17.47"
$ws.Range("E27").Value = "  -0.29%  "
$ws.Range("D28").Value = "0.06559"
$ws.Range("E28").Value = "  +15.16%  "
$ws.Range("E29").Value = "  +2.49%  "
$ws.Range("D30").Value = "1.487"
$ws.Range("E30").Value = "  +0.84%  "
$ws.Range("D31").Value = "4.074"
$ws.Range("E31").Value = "  -0.43%  "
$ws.Range("D32").Value = "4.052"
$ws.Range("E32").Value = "  +0.21%  "
$ws.Range("D33").Value = "1.841"
$ws.Range("E33").Value = "  +1.36%  "
$ws.Range("D34").Value = "1.142"
$ws.Range("E34").Value = "  -0.18%  "
$ws.Range("D35").Value = "0.6981"
$ws.Range("E35").Value = "  +1.02%  "
$ws.Range("D36").Value = "2.579"
$ws.Range("E36").Value = "  -0.23%  "
$ws.Range("D37").Value = "0.01860"
$ws.Range("E37").Value = "  +2.51%  "
$ws.Range("B38").Value = "MXToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D38").Value = "2.818"
$ws.Range("E38").Value = "  -0.71%  "
$ws.Range("B39").Value = "Maker"
$ws.Range("C39").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D39").Value = "1.248.44"
$ws.Range("E39").Value = "  -0.29%  "
$ws.Range("D40").Value = "6.781"
$ws.Range("D41").Value = "0.9349"
$ws.Range("E41").Value = "  +3.53%  "
$ws.Range("D42").Value = "1.000"
$ws.Range("E42").Value = "  +0.06%  "
$ws.Range("D43").Value = "2.015.18"
$ws.Range("E43").Value = "  +0.68%  "
$ws.Range("D44").Value = "101.11"
$ws.Range("E44").Value = "  -0.10%  "
$ws.Range("D45").Value = "65.43"
$ws.Range("E45").Value = "  -0.33%  "
$ws.Range("B46").Value = "Aptos"
$ws.Range("C46").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D46").Value = "7.065"
$ws.Range("E46").Value = "  +0.00%  "
$ws.Range("B47").Value = "BabyDogeCoin"
$ws.Range("C47").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D47").Value = "0.00000000118"
$ws.Range("E47").Value = "  +0.35%  "
$ws.Range("D48").Value = "1.716"
$ws.Range("E48").Value = "  +3.86%  "
$ws.Range("D49").Value = "9.041"
$ws.Range("E49").Value = "  +0.90%  "
$ws.Range("E50").Value = "  -1.13%  "
$ws.Range("D51").Value = "0.3913"
$ws.Range("E51").Value = "  -0.49%  "
